$d = $word.ActiveDocument

# Add a blank paragraph after the current last paragraph (bug #4). Using a
# bare carriage-return assigned to the collapsed end range (rather than
# InsertParagraphAfter) avoids leaving a stray empty run behind.
$p1 = $d.Paragraphs.Last.Range
$p1.Collapse(0)
$p1.Text = "`r"

# Add the new bug #5 paragraph after that blank paragraph.
$p2 = $d.Paragraphs.Last.Range
$p2.Collapse(0)
$p2.Text = "`r5. Despite the best efforts of trying (and multiple attempts at moving files around) the Dialog boxes have been incorrectly formatted. The css file is there, and has been moved around everywhere in an attempt to get it displaying correctly, but hasn't been resolved (as of yet)."

# Match the formatting of the surrounding paragraphs (Times New Roman,
# sz 24 / 12pt, incl. the complex-script variants so w:cs / w:szCs emit).
$finalRun = $d.Paragraphs.Last.Range
$finalRun.Font.Name = "Times New Roman"
$finalRun.Font.NameBi = "Times New Roman"
$finalRun.Font.Size = 12
$finalRun.Font.SizeBi = 12
